# Update the footer "date" placeholder text on the slide master and on
# every slide layout from 17/03/2023 -> 17/04/2023 (the datetimeFigureOut
# field was re-cached to a new date, as happens when PowerPoint refreshes
# the Header & Footer "Update automatically" field on save).

$p = $ppt.ActivePresentation
$newDate = "17/04/2023"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -ne $newDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master footer date placeholder.
Update-DatePlaceholders $p.SlideMaster.Shapes

# Every slide layout's own footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    $layout = $layouts.Item($l)
    Update-DatePlaceholders $layout.Shapes
}
